$ws = $excel.ActiveSheet

# Insert a new row at 214; this shifts the existing rows 214-272 down to 215-273.
$ws.Rows(214).Insert()

# Populate the newly inserted row 214 with the new weekly price record
# (Macroferia Regional de Talca - Poroto verde, origen Peru).
$ws.Range("A214").Value = 5
$ws.Range("B214").Value = "Macroferia Regional de Talca"
$ws.Range("C214").Value = "Maule"
$ws.Range("D214").Value = 45215
$ws.Range("E214").Value = 7
$ws.Range("F214").Value = 100112031
$ws.Range("G214").Value = "Poroto verde"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 200
$ws.Range("K214").Value = 30000
$ws.Range("L214").Value = 30000
$ws.Range("M214").Value = 30000
$ws.Range("N214").Value = '$/malla 25 kilos'
$ws.Range("O214").Value = "Perú"
$ws.Range("P214").Value = 1200
$ws.Range("Q214").Value = 25
$ws.Range("R214").Value = "Hortaliza"
